# Remove the "wbGetDataTask" Name/Value/Description rows from the Constants
# sheet. These occupied rows 23-24 (right after wbGetTransactionData and
# right before wbProcessTransaction). Deleting the two rows shifts
# everything below them up by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

$ws.Rows.Item(23).Resize(2).Delete() | Out-Null

# Leave the active selection where the author last left it.
$ws.Activate()
$ws.Range("B17").Select() | Out-Null
